# Apply weekly update: insert a new record row at row 57, pushing existing
# rows 57..153 down to 58..154 (preserving all their data), and populate the
# newly created row 57 with the new week's observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 57. This shifts rows 57-153 down
# to become rows 58-154 automatically (values + formatting).
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new data point.
$ws.Range("A57").Value = 7
$ws.Range("B57").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C57").Value = "Ñuble"
$ws.Range("D57").Value = 44533
$ws.Range("E57").Value = 16
$ws.Range("F57").Value = 100112017
$ws.Range("G57").Value = "Apio"
$ws.Range("H57").Value = "Americana (o)"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 60
$ws.Range("K57").Value = 8000
$ws.Range("L57").Value = 8500
$ws.Range("M57").Value = 8250
$ws.Range("N57").Value = "`$/docena de matas"
$ws.Range("O57").Value = "Provincia del Elquí"
$ws.Range("P57").Value = 1375
$ws.Range("Q57").Value = 6
$ws.Range("R57").Value = "Hortaliza"
